{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of `async (context) => { ... }`.\n//\n// The document's first paragraph is a centered date/weekday label; the\n// body also contains a single 20x5 table of two-digit multiplication\n// problems (\"A\u00d7B=C\"). This script:\n//   1. Replaces the date paragraph's text (keeping its run formatting).\n//   2. Replaces every table cell's text with its new value, preserving\n//      each cell's existing paragraph/run formatting (font, size, etc.)\n//      by assigning Cell.value rather than rebuilding the paragraph.\n\nconst NEW_DATE_TEXT = \"2023-04-06 Thursday\";\n\n// New cell text, in row-major order matching the table's 20 rows x 5 cols.\nconst NEW_TABLE_VALUES = [\n  [\n    \"46\u00d783=3818\",\n    \"84\u00d756=4704\",\n    \"82\u00d797=7954\",\n    \"86\u00d748=4128\",\n    \"72\u00d757=4104\"\n  ],\n  [\n    \"67\u00d737=2479\",\n    \"35\u00d743=1505\",\n    \"71\u00d795=6745\",\n    \"25\u00d798=2450\",\n    \"53\u00d722=1166\"\n  ],\n  [\n    \"53\u00d728=1484\",\n    \"75\u00d764=4800\",\n    \"24\u00d713=312\",\n    \"93\u00d710=930\",\n    \"34\u00d730=1020\"\n  ],\n  [\n    \"98\u00d762=6076\",\n    \"60\u00d741=2460\",\n    \"96\u00d756=5376\",\n    \"34\u00d740=1360\",\n    \"66\u00d765=4290\"\n  ],\n  [\n    \"31\u00d763=1953\",\n    \"35\u00d733=1155\",\n    \"51\u00d777=3927\",\n    \"61\u00d723=1403\",\n    \"78\u00d710=780\"\n  ],\n  [\n    \"53\u00d778=4134\",\n    \"44\u00d773=3212\",\n    \"44\u00d717=748\",\n    \"32\u00d775=2400\",\n    \"86\u00d794=8084\"\n  ],\n  [\n    \"21\u00d794=1974\",\n    \"82\u00d792=7544\",\n    \"74\u00d784=6216\",\n    \"54\u00d778=4212\",\n    \"78\u00d795=7410\"\n  ],\n  [\n    \"40\u00d711=440\",\n    \"37\u00d789=3293\",\n    \"53\u00d760=3180\",\n    \"15\u00d735=525\",\n    \"74\u00d710=740\"\n  ],\n  [\n    \"27\u00d729=783\",\n    \"95\u00d726=2470\",\n    \"46\u00d712=552\",\n    \"18\u00d736=648\",\n    \"91\u00d734=3094\"\n  ],\n  [\n    \"79\u00d7100=7900\",\n    \"18\u00d726=468\",\n    \"14\u00d785=1190\",\n    \"78\u00d743=3354\",\n    \"82\u00d739=3198\"\n  ],\n  [\n    \"19\u00d743=817\",\n    \"35\u00d724=840\",\n    \"32\u00d744=1408\",\n    \"64\u00d749=3136\",\n    \"37\u00d782=3034\"\n  ],\n  [\n    \"68\u00d742=2856\",\n    \"38\u00d779=3002\",\n    \"25\u00d792=2300\",\n    \"64\u00d788=5632\",\n    \"15\u00d775=1125\"\n  ],\n  [\n    \"94\u00d729=2726\",\n    \"98\u00d728=2744\",\n    \"30\u00d798=2940\",\n    \"72\u00d726=1872\",\n    \"38\u00d752=1976\"\n  ],\n  [\n    \"48\u00d776=3648\",\n    \"29\u00d746=1334\",\n    \"55\u00d763=3465\",\n    \"13\u00d715=195\",\n    \"98\u00d742=4116\"\n  ],\n  [\n    \"95\u00d774=7030\",\n    \"38\u00d728=1064\",\n    \"52\u00d711=572\",\n    \"47\u00d769=3243\",\n    \"60\u00d733=1980\"\n  ],\n  [\n    \"95\u00d742=3990\",\n    \"13\u00d778=1014\",\n    \"27\u00d798=2646\",\n    \"100\u00d736=3600\",\n    \"88\u00d778=6864\"\n  ],\n  [\n    \"73\u00d722=1606\",\n    \"78\u00d751=3978\",\n    \"13\u00d784=1092\",\n    \"21\u00d768=1428\",\n    \"79\u00d745=3555\"\n  ],\n  [\n    \"17\u00d729=493\",\n    \"71\u00d760=4260\",\n    \"20\u00d730=600\",\n    \"24\u00d719=456\",\n    \"36\u00d798=3528\"\n  ],\n  [\n    \"27\u00d795=2565\",\n    \"95\u00d731=2945\",\n    \"52\u00d783=4316\",\n    \"15\u00d797=1455\",\n    \"63\u00d797=6111\"\n  ],\n  [\n    \"95\u00d745=4275\",\n    \"58\u00d775=4350\",\n    \"45\u00d768=3060\",\n    \"70\u00d795=6650\",\n    \"94\u00d756=5264\"\n  ]\n];\n\nconst body = context.document.body;\n\n// --- 1. Update the date/weekday paragraph (first paragraph in the body) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(NEW_DATE_TEXT, Word.InsertLocation.replace);\n\n// --- 2. Update every cell in the multiplication table ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst rowCount = NEW_TABLE_VALUES.length;\nconst columnCount = NEW_TABLE_VALUES[0].length;\n\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = NEW_TABLE_VALUES[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $d (= $word.ActiveDocument) is the document already opened for us.\n#\n# The document body has a single centered paragraph with a date/weekday\n# label, followed by one 20-row x 5-column table of two-digit\n# multiplication problems (\"A\u00d7B=C\"). This script replaces the date text\n# and every table cell's text with new values, using Range.Text\n# assignment so each cell/paragraph keeps its existing run formatting\n# (font, size, etc.).\n\n$d = $word.ActiveDocument\n\n# New text for the date/weekday paragraph (first paragraph in the body).\n$newDateText = \"2023-04-06 Thursday\"\n\n# New cell text, in row-major order matching the table's 20 rows x 5 cols.\n$newTableValues = @(\n    @(\"46\u00d783=3818\", \"84\u00d756=4704\", \"82\u00d797=7954\", \"86\u00d748=4128\", \"72\u00d757=4104\"),\n    @(\"67\u00d737=2479\", \"35\u00d743=1505\", \"71\u00d795=6745\", \"25\u00d798=2450\", \"53\u00d722=1166\"),\n    @(\"53\u00d728=1484\", \"75\u00d764=4800\", \"24\u00d713=312\", \"93\u00d710=930\", \"34\u00d730=1020\"),\n    @(\"98\u00d762=6076\", \"60\u00d741=2460\", \"96\u00d756=5376\", \"34\u00d740=1360\", \"66\u00d765=4290\"),\n    @(\"31\u00d763=1953\", \"35\u00d733=1155\", \"51\u00d777=3927\", \"61\u00d723=1403\", \"78\u00d710=780\"),\n    @(\"53\u00d778=4134\", \"44\u00d773=3212\", \"44\u00d717=748\", \"32\u00d775=2400\", \"86\u00d794=8084\"),\n    @(\"21\u00d794=1974\", \"82\u00d792=7544\", \"74\u00d784=6216\", \"54\u00d778=4212\", \"78\u00d795=7410\"),\n    @(\"40\u00d711=440\", \"37\u00d789=3293\", \"53\u00d760=3180\", \"15\u00d735=525\", \"74\u00d710=740\"),\n    @(\"27\u00d729=783\", \"95\u00d726=2470\", \"46\u00d712=552\", \"18\u00d736=648\", \"91\u00d734=3094\"),\n    @(\"79\u00d7100=7900\", \"18\u00d726=468\", \"14\u00d785=1190\", \"78\u00d743=3354\", \"82\u00d739=3198\"),\n    @(\"19\u00d743=817\", \"35\u00d724=840\", \"32\u00d744=1408\", \"64\u00d749=3136\", \"37\u00d782=3034\"),\n    @(\"68\u00d742=2856\", \"38\u00d779=3002\", \"25\u00d792=2300\", \"64\u00d788=5632\", \"15\u00d775=1125\"),\n    @(\"94\u00d729=2726\", \"98\u00d728=2744\", \"30\u00d798=2940\", \"72\u00d726=1872\", \"38\u00d752=1976\"),\n    @(\"48\u00d776=3648\", \"29\u00d746=1334\", \"55\u00d763=3465\", \"13\u00d715=195\", \"98\u00d742=4116\"),\n    @(\"95\u00d774=7030\", \"38\u00d728=1064\", \"52\u00d711=572\", \"47\u00d769=3243\", \"60\u00d733=1980\"),\n    @(\"95\u00d742=3990\", \"13\u00d778=1014\", \"27\u00d798=2646\", \"100\u00d736=3600\", \"88\u00d778=6864\"),\n    @(\"73\u00d722=1606\", \"78\u00d751=3978\", \"13\u00d784=1092\", \"21\u00d768=1428\", \"79\u00d745=3555\"),\n    @(\"17\u00d729=493\", \"71\u00d760=4260\", \"20\u00d730=600\", \"24\u00d719=456\", \"36\u00d798=3528\"),\n    @(\"27\u00d795=2565\", \"95\u00d731=2945\", \"52\u00d783=4316\", \"15\u00d797=1455\", \"63\u00d797=6111\"),\n    @(\"95\u00d745=4275\", \"58\u00d775=4350\", \"45\u00d768=3060\", \"70\u00d795=6650\", \"94\u00d756=5264\")\n)\n\n# --- 1. Update the date/weekday paragraph ---\n$d.Paragraphs.Item(1).Range.Text = $newDateText\n\n# --- 2. Update every cell in the multiplication table ---\n$table = $d.Tables.Item(1)\n$rowCount = $newTableValues.Count\nfor ($r = 1; $r -le $rowCount; $r++) {\n    $rowValues = $newTableValues[$r - 1]\n    $colCount = $rowValues.Count\n    for ($c = 1; $c -le $colCount; $c++) {\n        $table.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
